# Atualização automática de CAMPINA_DAS_MISSOES.xlsx
#
# 1. Remove the "Desarquivamentos Pendentes" sheet (no longer needed).
# 2. Rename "Paineis DARQ" -> "PAINEIS DARQ".
# 3. Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO".

$wb = $excel.ActiveWorkbook

# 1. Delete the obsolete "Desarquivamentos Pendentes" worksheet.
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null

# 2. Rename "Paineis DARQ" to all caps "PAINEIS DARQ".
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# 3. Rename "Recolhimento x Eliminacao" to all caps (with accent) "RECOLHIMENTO X ELIMINAÇÃO".
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
